# "Samenvatting notulen interview toegevoegd"
#
# - Strike through the "Wasmachine stuurt errorrapport naar SWIRL wat er
#   fout is gegaan." bullet (text + its paragraph mark).
# - Word's auto-maintained "_GoBack" bookmark (last-edit location) moves
#   from the old edit spot ("Samenvatting nog mailen.") to now wrap the
#   paragraph that was just edited.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Wasmachine stuurt errorrapport naar SWIRL*") {
        $target = $p
    }
}

# Strike through the run text and the paragraph mark itself.
$target.Range.Font.StrikeThrough = $true

# Re-seat "_GoBack" (it already exists elsewhere in the document) onto the
# paragraph that was just edited - it covers the whole paragraph, mark
# included, which is why its end lands at the very start of the next one.
$d.Bookmarks.Add("_GoBack", $target.Range)
